# Auto-generated: update computed market-price columns (H-N) per the
# scheduled Sheets refresh. Values come from the commit's canonical OOXML.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2129154.2
$ws.Range("J17").Value = 2129154.2
$ws.Range("L17").Value = 6387462.600000001
$ws.Range("N17").Value = -6387798.600000001
# Row 19
$ws.Range("H19").Value = 382.16666
$ws.Range("I19").Value = 339.5
$ws.Range("K19").Value = 339.5
$ws.Range("M19").Value = -164.5
# Row 69
$ws.Range("H69").Value = 8344584.5
$ws.Range("I69").Value = 16674194
$ws.Range("J69").Value = 14974.833
$ws.Range("K69").Value = 50022582
$ws.Range("L69").Value = 44924.499
$ws.Range("M69").Value = -50021708
$ws.Range("N69").Value = -46672.499
# Row 72
$ws.Range("H72").Value = 8344584.5
$ws.Range("I72").Value = 16674194
$ws.Range("J72").Value = 14974.833
$ws.Range("K72").Value = 150067746
$ws.Range("L72").Value = 134773.497
$ws.Range("M72").Value = -150063378
$ws.Range("N72").Value = -143509.497
# Row 113
$ws.Range("H113").Value = 8023.8335
$ws.Range("I113").Value = 7277.273
$ws.Range("K113").Value = 7277.273
$ws.Range("M113").Value = -4023.273
# Row 118
$ws.Range("H118").Value = 787.8125
$ws.Range("I118").Value = 662
$ws.Range("J118").Value = 1333
$ws.Range("K118").Value = 1986
$ws.Range("L118").Value = 3999
$ws.Range("M118").Value = -329
$ws.Range("N118").Value = -7313
# Row 135
$ws.Range("H135").Value = 2337.5
$ws.Range("I135").Value = 1764.909
$ws.Range("K135").Value = 15884.181
$ws.Range("M135").Value = -13349.181
# Row 137
$ws.Range("H137").Value = 7144844.5
$ws.Range("I137").Value = 2140.2307
$ws.Range("K137").Value = 6420.6921
$ws.Range("M137").Value = -3870.6921
# Row 138
$ws.Range("H138").Value = 12198.919
$ws.Range("I138").Value = 16024.2
$ws.Range("J138").Value = 11769.112
$ws.Range("K138").Value = 48072.60000000001
$ws.Range("L138").Value = 35307.336
$ws.Range("M138").Value = -42932.60000000001
$ws.Range("N138").Value = -45587.336
# Row 141
$ws.Range("H141").Value = 977.75
$ws.Range("I141").Value = 978.54285
$ws.Range("J141").Value = 950
$ws.Range("K141").Value = 2935.62855
$ws.Range("L141").Value = 2850
$ws.Range("M141").Value = 2244.37145
$ws.Range("N141").Value = -13210

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 511.22223
$ws.Range("I4").Value = 511.22223
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 511.22223
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -395.22223
$ws.Range("N4").ClearContents()
# Row 32
$ws.Range("H32").Value = 4598164.5
$ws.Range("I32").Value = 5929061
$ws.Range("J32").Value = 57458.707
$ws.Range("K32").Value = 5929061
$ws.Range("L32").Value = 57458.707
$ws.Range("M32").Value = -5928774
$ws.Range("N32").Value = -58032.707
# Row 61
$ws.Range("H61").Value = 1581896.4
$ws.Range("I61").Value = 8060.591
$ws.Range("K61").Value = 8060.591
$ws.Range("M61").Value = -7848.591
# Row 74
$ws.Range("H74").Value = 305079.03
$ws.Range("I74").Value = 4063.3728
$ws.Range("K74").Value = 4063.3728
$ws.Range("M74").Value = -3189.3728
# Row 77
$ws.Range("H77").Value = 305079.03
$ws.Range("I77").Value = 4063.3728
$ws.Range("K77").Value = 20316.864
$ws.Range("M77").Value = -15948.864
# Row 102
$ws.Range("H102").Value = 47622540
$ws.Range("I102").Value = 66671036
$ws.Range("K102").Value = 66671036
$ws.Range("M102").Value = -66669414
# Row 132
$ws.Range("H132").Value = 2205.3586
$ws.Range("I132").Value = 1831.2307
$ws.Range("K132").Value = 5493.6921
$ws.Range("M132").Value = -2963.6921
# Row 136
$ws.Range("H136").Value = 1581896.4
$ws.Range("I136").Value = 8060.591
$ws.Range("K136").Value = 24181.773
$ws.Range("M136").Value = -21631.773

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 47425560
$ws.Range("I134").Value = 71174.07000000001
$ws.Range("J134").Value = 225004500
$ws.Range("K134").Value = 213522.21
$ws.Range("L134").Value = 675013500
$ws.Range("M134").Value = -210987.21
$ws.Range("N134").Value = -675018570

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2610.4736
$ws.Range("I22").Value = 533.2778
$ws.Range("K22").Value = 533.2778
$ws.Range("M22").Value = -183.2778
# Row 31
$ws.Range("H31").Value = 4343.087
$ws.Range("I31").Value = 3816.647
$ws.Range("J31").Value = 5834.6665
$ws.Range("K31").Value = 3816.647
$ws.Range("L31").Value = 5834.6665
$ws.Range("M31").Value = -3521.647
$ws.Range("N31").Value = -6424.6665
# Row 34
$ws.Range("H34").Value = 4343.087
$ws.Range("I34").Value = 3816.647
$ws.Range("J34").Value = 5834.6665
$ws.Range("K34").Value = 3816.647
$ws.Range("L34").Value = 5834.6665
$ws.Range("M34").Value = -3614.647
$ws.Range("N34").Value = -6238.6665
# Row 35
$ws.Range("H35").Value = 2023.2916
$ws.Range("J35").Value = 1403
$ws.Range("L35").Value = 1403
$ws.Range("N35").Value = -1991
# Row 58
$ws.Range("H58").Value = 4166.227
$ws.Range("I58").Value = 1612.4
$ws.Range("J58").Value = 4917.353
$ws.Range("K58").Value = 1612.4
$ws.Range("L58").Value = 4917.353
$ws.Range("M58").Value = -1409.4
$ws.Range("N58").Value = -5323.353
# Row 134
$ws.Range("H134").Value = 2756.157
$ws.Range("I134").Value = 1889.079
$ws.Range("J134").Value = 5290.6924
$ws.Range("K134").Value = 5667.237
$ws.Range("L134").Value = 15872.0772
$ws.Range("M134").Value = -3132.237
$ws.Range("N134").Value = -20942.0772
# Row 136
$ws.Range("H136").Value = 4166.227
$ws.Range("I136").Value = 1612.4
$ws.Range("J136").Value = 4917.353
$ws.Range("K136").Value = 4837.200000000001
$ws.Range("L136").Value = 14752.059
$ws.Range("M136").Value = -2287.200000000001
$ws.Range("N136").Value = -19852.059
# Row 141
$ws.Range("H141").Value = 539755.25
$ws.Range("J141").Value = 564506.2
$ws.Range("L141").Value = 564506.2
$ws.Range("N141").Value = -574866.2

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 29007600
$ws.Range("I131").Value = 22731772
$ws.Range("J131").Value = 37375372
$ws.Range("K131").Value = 68195316
$ws.Range("L131").Value = 112126116
$ws.Range("M131").Value = -68190276
$ws.Range("N131").Value = -112136196
# Row 132
$ws.Range("H132").Value = 1991.1818
$ws.Range("J132").Value = 1990.3
$ws.Range("L132").Value = 17912.7
$ws.Range("N132").Value = -22972.7

$ws = $wb.Worksheets.Item("GSM")
# Row 99
$ws.Range("H99").Value = 23535.777
$ws.Range("I99").Value = 20260.428
$ws.Range("K99").Value = 20260.428
$ws.Range("M99").Value = -18014.428
# Row 132
$ws.Range("H132").Value = 6293819.5
$ws.Range("I132").Value = 7168.231
$ws.Range("J132").Value = 15374538
$ws.Range("K132").Value = 21504.693
$ws.Range("L132").Value = 46123614
$ws.Range("M132").Value = -18974.693
$ws.Range("N132").Value = -46128674

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 10486.308
$ws.Range("I46").Value = 11633.4
$ws.Range("K46").Value = 11633.4
$ws.Range("M46").Value = -11445.4
# Row 55
$ws.Range("H55").Value = 767.4375
$ws.Range("I55").Value = 251.28572
$ws.Range("J55").Value = 1168.8889
$ws.Range("K55").Value = 251.28572
$ws.Range("L55").Value = 1168.8889
$ws.Range("M55").Value = -78.28572
$ws.Range("N55").Value = -1514.8889
# Row 61
$ws.Range("H61").Value = 1698.7142
$ws.Range("I61").Value = 1581.6097
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 1581.6097
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -1379.6097
$ws.Range("N61").Value = -6904
# Row 82
$ws.Range("H82").Value = 7162.6665
$ws.Range("I82").Value = 1332.3334
$ws.Range("J82").Value = 10077.833
$ws.Range("K82").Value = 1332.3334
$ws.Range("L82").Value = 10077.833
$ws.Range("M82").Value = -971.3334
$ws.Range("N82").Value = -10799.833
# Row 85
$ws.Range("H85").Value = 7162.6665
$ws.Range("I85").Value = 1332.3334
$ws.Range("J85").Value = 10077.833
$ws.Range("K85").Value = 1332.3334
$ws.Range("L85").Value = 10077.833
$ws.Range("M85").Value = -84.33339999999998
$ws.Range("N85").Value = -12573.833
# Row 113
$ws.Range("H113").Value = 1698.7142
$ws.Range("I113").Value = 1581.6097
$ws.Range("J113").Value = 6500
$ws.Range("K113").Value = 1581.6097
$ws.Range("L113").Value = 6500
$ws.Range("M113").Value = 588.3903
$ws.Range("N113").Value = -10840
# Row 136
$ws.Range("H136").Value = 7687.5713
$ws.Range("I136").Value = 8211.666999999999
$ws.Range("K136").Value = 24635.001
$ws.Range("M136").Value = -22085.001

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 6466.5
$ws.Range("I62").Value = 5999.6665
$ws.Range("K62").Value = 5999.6665
$ws.Range("M62").Value = -5375.6665
# Row 65
$ws.Range("H65").Value = 6466.5
$ws.Range("I65").Value = 5999.6665
$ws.Range("K65").Value = 29998.3325
$ws.Range("M65").Value = -26878.3325
# Row 81
$ws.Range("H81").Value = 12788.777
$ws.Range("J81").Value = 2600
$ws.Range("L81").Value = 5200
$ws.Range("N81").Value = -7322
# Row 84
$ws.Range("H84").Value = 12788.777
$ws.Range("J84").Value = 2600
$ws.Range("L84").Value = 26000
$ws.Range("N84").Value = -36608
# Row 122
$ws.Range("H122").Value = 7182.4614
$ws.Range("I122").Value = 7182.4614
$ws.Range("K122").Value = 21547.3842
$ws.Range("M122").Value = -19097.3842
# Row 132
$ws.Range("H132").Value = 43373.543
$ws.Range("I132").Value = 144359.78
$ws.Range("J132").Value = 1790.9706
$ws.Range("K132").Value = 433079.34
$ws.Range("L132").Value = 5372.9118
$ws.Range("M132").Value = -430549.34
$ws.Range("N132").Value = -10432.9118
# Row 136
$ws.Range("H136").Value = 20334.246
$ws.Range("I136").Value = 20738.037
$ws.Range("K136").Value = 62214.111
$ws.Range("M136").Value = -59664.111

